# Updated symbol list on Thu Dec 29 18:41:11 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Keep a reference to an untouched "normal" (default) cell style so that
# after forcing text storage for numeric-looking values we can restore the
# original (unstyled) look of each cell.
$normalStyle = $ws.Range("D3").Style

# Price ("D" column) values are stored as text in the workbook (inlineStr),
# so updating them with plain numeric-looking strings would make Excel
# silently reinterpret them as numbers. Force each target cell to Text
# format first, assign the new text value, then restore its original style
# so no unintended formatting change is left behind.
$priceUpdates = [ordered]@{
    "D2"  = "245.73"
    "D4"  = "5.278"
    "D5"  = "0.05774"
    "D6"  = "6.504"
    "D7"  = "3.144"
    "D8"  = "0.8114"
    "D9"  = "0.8622"
    "D11" = "0.06953"
    "D12" = "0.03146"
    "D13" = "0.02914"
    "D14" = "0.09378"
    "D15" = "3.757"
    "D16" = "0.001531"
    "D17" = "0.04694"
    "D18" = "0.0006012"
    "D19" = "0.006142"
    "D21" = "0.004640"
    "D23" = "3.502"
    "D24" = "2.149"
    "D28" = "0.0002332"
    "D40" = "0.03709"
    "D41" = "0.006383"
    "D42" = "0.1054"
    "D43" = "0.003001"
    "D44" = "0.007756"
    "D45" = "0.00005253"
    "D47" = "0.4402"
    "D48" = "0.002483"
}

foreach ($cellRef in $priceUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$cellRef]
    $cell.Style = $normalStyle
}

# Volume(1h) ("E" column) text labels that changed wording.
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E22").Value = "21NitroExNTXWorstin24h"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("E43").Value = "42CEJICEJI"
